$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.837.90"
$ws.Range("E2").Value = '  +5.69%  '
$ws.Range("D3").Value = "'2.592.22"
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").Value = "'591.48"
$ws.Range("E5").Value = '  +4.12%  '
$ws.Range("D6").Value = "'156.00"
$ws.Range("E6").Value = '  +7.51%  '
$ws.Range("D7").Value = "'0.997"
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("D8").Value = "'0.548"
$ws.Range("E8").Value = '  +3.87%  '
$ws.Range("D9").Value = "'2.623.09"
$ws.Range("E9").Value = '  +7.45%  '
$ws.Range("E10").Value = '  +5.49%  '
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("E12").Value = '  +5.41%  '
$ws.Range("D13").Value = "'5.34"
$ws.Range("E13").Value = '  +3.01%  '
$ws.Range("D14").Value = "'29.48"
$ws.Range("E14").Value = '  +3.61%  '
$ws.Range("D15").Value = "'0.0000182"
$ws.Range("E15").Value = '  +5.78%  '
$ws.Range("D16").Value = "'3.063.15"
$ws.Range("E16").Value = '  +6.15%  '
$ws.Range("D17").Value = "'65.666.95"
$ws.Range("E17").Value = '  +5.70%  '
$ws.Range("D18").Value = "'2.618.63"
$ws.Range("E18").Value = '  +6.96%  '
$ws.Range("D19").Value = "'8.26"
$ws.Range("E19").Value = '  +8.47%  '
$ws.Range("D20").Value = "'11.27"
$ws.Range("E20").Value = '  +5.77%  '
$ws.Range("D21").Value = "'355.02"
$ws.Range("E21").Value = '  +11.21%  '
$ws.Range("E22").Value = '  +5.18%  '
$ws.Range("E23").Value = '  +4.59%  '
$ws.Range("D24").Value = "'0.998"
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").Value = "'10.03"
$ws.Range("E25").Value = '  +3.30%  '
$ws.Range("D26").Value = "'66.29"
$ws.Range("E26").Value = '  +2.37%  '
$ws.Range("D27").Value = "'636.65"
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").Value = "'0.0000106"
$ws.Range("E28").Value = '  +12.40%  '
$ws.Range("D29").Value = "'2.723.32"
$ws.Range("E29").Value = '  +6.37%  '
$ws.Range("D30").Value = "'1.51"
$ws.Range("E30").Value = '  +8.01%  '
$ws.Range("D31").Value = "'0.991"
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("D32").Value = "'8.25"
$ws.Range("E32").Value = '  +5.62%  '
$ws.Range("E33").Value = '  +6.51%  '
$ws.Range("D34").Value = "'0.138"
$ws.Range("E34").Value = '  +5.25%  '
$ws.Range("E35").Value = '  +10.76%  '
$ws.Range("D36").Value = "'0.995"
$ws.Range("E36").Value = '  -0.31%  '
$ws.Range("E37").Value = '  +8.47%  '
$ws.Range("D38").Value = "'5.67"
$ws.Range("E38").Value = '  +8.29%  '
$ws.Range("D39").Value = "'2.91"
$ws.Range("E39").Value = '  +9.04%  '
$ws.Range("D40").Value = "'19.43"
$ws.Range("E40").Value = '  +5.68%  '
$ws.Range("D41").Value = "'156.00"
$ws.Range("E41").Value = '  +3.73%  '
$ws.Range("E42").Value = '  +3.56%  '
$ws.Range("E43").Value = '  +7.62%  '
$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").Value = "'0.0₆0323"
$ws.Range("E44").Value = '  +5.68%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = "'42.26"
$ws.Range("E45").Value = '  +1.23%  '
$ws.Range("D46").Value = "'163.33"
$ws.Range("E46").Value = '  +8.15%  '
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("D48").Value = "'16.19"
$ws.Range("E48").Value = '  +5.66%  '
$ws.Range("D49").Value = "'3.77"
$ws.Range("E49").Value = '  +7.57%  '
$ws.Range("D50").Value = "'21.83"
$ws.Range("E50").Value = '  +9.65%  '
$ws.Range("D51").Value = "'0.640"
$ws.Range("E51").Value = '  +6.64%  '
